$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.215.59'
$ws.Range('E2').Value = '  +3.87%  '

$ws.Range('D3').Value = '2.433.11'
$ws.Range('E3').Value = '  +5.16%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = "'556.68"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.78%  '

$ws.Range('D6').Value = "'138.87"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.78%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('E8').Value = '  +1.19%  '

$ws.Range('D9').Value = '2.430.17'
$ws.Range('E9').Value = '  +5.11%  '

$ws.Range('E10').Value = '  +3.54%  '

$ws.Range('E11').Value = '  +4.40%  '

$ws.Range('E12').Value = '  +0.32%  '

$ws.Range('E13').Value = '  +5.03%  '

$ws.Range('D14').Value = "'26.15"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +11.96%  '

$ws.Range('D15').Value = '2.868.49'
$ws.Range('E15').Value = '  +5.22%  '

$ws.Range('D16').Value = '62.075.06'
$ws.Range('E16').Value = '  +3.67%  '

$ws.Range('E17').Value = '  +8.01%  '

$ws.Range('D18').Value = '2.431.51'
$ws.Range('E18').Value = '  +4.08%  '

$ws.Range('D19').Value = "'11.21"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.08%  '

$ws.Range('D20').Value = "'346.06"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +10.91%  '

$ws.Range('E21').Value = '  +3.07%  '

$ws.Range('D22').Value = "'6.78"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.36%  '

$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('D24').Value = "'65.14"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.39%  '

$ws.Range('E25').Value = '  +1.32%  '

$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.14%  '

$ws.Range('E27').Value = '  +14.22%  '

$ws.Range('E28').Value = '  +5.98%  '

$ws.Range('D29').Value = "'1.33"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +14.37%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = "'1.80"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.96%  '

$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0783'
$ws.Range('E31').Value = '  +8.49%  '

$ws.Range('E32').Value = '  +9.63%  '

$ws.Range('D33').Value = "'171.76"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.43%  '

$ws.Range('E34').Value = '  +5.52%  '

$ws.Range('E35').Value = '  +4.96%  '

$ws.Range('D36').Value = "'18.57"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.09%  '

$ws.Range('E37').Value = '  +11.93%  '

$ws.Range('D38').Value = "'365.83"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +15.63%  '

$ws.Range('E40').Value = '  -0.04%  '

$ws.Range('E41').Value = '  +11.63%  '

$ws.Range('D42').Value = "'39.21"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.66%  '

$ws.Range('D43').Value = "'146.26"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.29%  '

$ws.Range('D44').Value = "'3.67"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.04%  '

$ws.Range('D45').Value = "'20.55"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +9.95%  '

$ws.Range('E46').Value = '  +1.73%  '

$ws.Range('E47').Value = '  +4.82%  '

$ws.Range('E48').Value = '  +5.70%  '

$ws.Range('E49').Value = '  +4.94%  '

$ws.Range('D50').Value = "'17.87"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.98%  '

$ws.Range('E51').Value = '  -2.20%  '
